# Update '想去人数' (wanted-attendance) counters for all four sheets
# to the freshly scraped values (gh-pages data refresh at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 2217
$ws.Range("F5").Value = 4365
$ws.Range("F8").Value = 1358
$ws.Range("F12").Value = 412
$ws.Range("F13").Value = 680131
$ws.Range("F14").Value = 1682
$ws.Range("F15").Value = 592
$ws.Range("F16").Value = 1494
$ws.Range("F20").Value = 2325
$ws.Range("F21").Value = 1166
$ws.Range("F22").Value = 2736
$ws.Range("F23").Value = 1584
$ws.Range("F24").Value = 905
$ws.Range("F25").Value = 1590
$ws.Range("F26").Value = 543
$ws.Range("F27").Value = 1102
$ws.Range("F28").Value = 1136
$ws.Range("F29").Value = 1109
$ws.Range("F32").Value = 2063
$ws.Range("F33").Value = 591
$ws.Range("F34").Value = 1351
$ws.Range("F35").Value = 3164
$ws.Range("F37").Value = 1157
$ws.Range("F40").Value = 2652
$ws.Range("F42").Value = 1005
$ws.Range("F44").Value = 1030
$ws.Range("F46").Value = 885
$ws.Range("F47").Value = 169
$ws.Range("F48").Value = 675

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 145561
$ws.Range("F10").Value = 145561
$ws.Range("F16").Value = 236
$ws.Range("F17").Value = 344
$ws.Range("F19").Value = 428
$ws.Range("F20").Value = 186
$ws.Range("F24").Value = 674
$ws.Range("F25").Value = 93
$ws.Range("F29").Value = 383
$ws.Range("F30").Value = 286
$ws.Range("F32").Value = 77
$ws.Range("F33").Value = 77

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 3142
$ws.Range("F5").Value = 260
$ws.Range("F6").Value = 14
$ws.Range("F7").Value = 839
$ws.Range("F8").Value = 1236
$ws.Range("F9").Value = 652
$ws.Range("F10").Value = 1624
$ws.Range("F11").Value = 157
$ws.Range("F12").Value = 2064

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 839
$ws.Range("F3").Value = 652
$ws.Range("F5").Value = 1624
$ws.Range("F7").Value = 2217
$ws.Range("F8").Value = 157
$ws.Range("F9").Value = 2064
$ws.Range("F10").Value = 4365
$ws.Range("F12").Value = 1358
$ws.Range("F15").Value = 412
$ws.Range("F16").Value = 680134
$ws.Range("F19").Value = 1682
$ws.Range("F20").Value = 145561
$ws.Range("F21").Value = 1494
$ws.Range("F25").Value = 2325
$ws.Range("F26").Value = 1166
$ws.Range("F27").Value = 2736
$ws.Range("F28").Value = 1584
$ws.Range("F29").Value = 905
$ws.Range("F31").Value = 1590
$ws.Range("F32").Value = 543
$ws.Range("F33").Value = 186
$ws.Range("F34").Value = 1102
$ws.Range("F35").Value = 1136
$ws.Range("F36").Value = 1109
$ws.Range("F38").Value = 2063
$ws.Range("F39").Value = 1351
$ws.Range("F40").Value = 3164
$ws.Range("F42").Value = 1157
$ws.Range("F43").Value = 383
$ws.Range("F44").Value = 286
$ws.Range("F45").Value = 77
$ws.Range("F46").Value = 2652
$ws.Range("F48").Value = 1005
$ws.Range("F51").Value = 1030
$ws.Range("F52").Value = 169
$ws.Range("F53").Value = 675

Write-Output "Updated 86 cells across 4 sheets."
